{"js": "// Wordage fix (Ryan's 6/27/2023 email): correct the 90% confidence interval\n// for the Oreo filling weight problem from (2.808, 2.988) to (2.535, 3.165).\nconst body = context.document.body;\n\n// Replace every occurrence of the old lower bound with the new one.\nconst lowerBound = body.search(\"2.808\", { matchCase: true });\nlowerBound.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < lowerBound.items.length; i++) {\n  lowerBound.items[i].insertText(\"2.535\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Replace every occurrence of the old upper bound with the new one.\nconst upperBound = body.search(\"2.988\", { matchCase: true });\nupperBound.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < upperBound.items.length; i++) {\n  upperBound.items[i].insertText(\"3.165\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Wordage fix (Ryan's 6/27/2023 email): correct the 90% confidence interval\n# for the Oreo filling weight problem from (2.808, 2.988) to (2.535, 3.165).\n$d = $word.ActiveDocument\n\n# Replace every occurrence of the old lower bound with the new one.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n  \"2.808\",          # FindText\n  $false,           # MatchCase\n  $false,           # MatchWholeWord\n  $false,           # MatchWildcards\n  $false,           # MatchSoundsLike\n  $false,           # MatchAllWordForms\n  $true,            # Forward\n  1,                # Wrap (wdFindContinue)\n  $false,           # Format\n  \"2.535\",          # ReplaceWith\n  2                 # Replace (wdReplaceAll)\n)\n\n# Replace every occurrence of the old upper bound with the new one.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n  \"2.988\",          # FindText\n  $false,           # MatchCase\n  $false,           # MatchWholeWord\n  $false,           # MatchWildcards\n  $false,           # MatchSoundsLike\n  $false,           # MatchAllWordForms\n  $true,            # Forward\n  1,                # Wrap (wdFindContinue)\n  $false,           # Format\n  \"3.165\",          # ReplaceWith\n  2                 # Replace (wdReplaceAll)\n)\n"}
